# KULR & IBM Update
# Applies the Q1'24 ("M" column) actuals + Q3'24 ("Q" column) revisions to the
# "Model" sheet, and rolls the RPO/cash/debt reconciliation on "Main" forward.

$wb = $excel.ActiveWorkbook

$wsMain  = $wb.Worksheets.Item("Main")
$wsModel = $wb.Worksheets.Item("Model")

# ---------------------------------------------------------------------------
# Main sheet: RPO roll-forward
# ---------------------------------------------------------------------------
$wsMain.Range("N5").Formula = "=13197+17+505"
$wsMain.Range("N6").Formula = "=52980+3599"

# O3/O4/O5/O6 pick up the right-aligned "unit label" style used elsewhere.
$wsMain.Range("O3").HorizontalAlignment = -4152
$wsMain.Range("O4").HorizontalAlignment = -4152
$wsMain.Range("O5").HorizontalAlignment = -4152
$wsMain.Range("O6").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Model sheet: fill in Q1'24 actuals (column M) for the segment detail that
# was previously left blank, and true-up a couple of Q3'24 (column Q) figures.
# ---------------------------------------------------------------------------

# Segment revenue detail (rows 5-13) - column M (Q1'24) and Q (Q3'24).
$wsModel.Range("M5").Value  = 4187
$wsModel.Range("Q5").Value  = 4600
$wsModel.Range("M6").Value  = 1759
$wsModel.Range("Q6").Value  = 1925
$wsModel.Range("M7").Value  = 2291
$wsModel.Range("Q7").Value  = 2327
$wsModel.Range("M8").Value  = 1944
$wsModel.Range("Q8").Value  = 1921
$wsModel.Range("M9").Value  = 943
$wsModel.Range("Q9").Value  = 905
$wsModel.Range("M10").Value = 1943
$wsModel.Range("Q10").Value = 1765
$wsModel.Range("M11").Value = 1329
$wsModel.Range("Q11").Value = 1277
$wsModel.Range("M12").Value = 186
$wsModel.Range("Q12").Value = 181
$wsModel.Range("M13").Value = 170
$wsModel.Range("Q13").Value = 68

# Helper: copy the number format + bold-ness from the "P" column (the most
# recently-populated quarter that uses the same style family as the new "M"
# column) so the new cells land on the existing style ids instead of minting
# duplicates.
function Copy-ModelStyle($row) {
    $src = $wsModel.Range("P$row")
    $dst = $wsModel.Range("M$row")
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Bold = $src.Font.Bold
    # Only force an explicit alignment when the source actually has one -
    # writing back "General" (the default) mints a spurious new style.
    if ($src.HorizontalAlignment -eq -4152) {
        $dst.HorizontalAlignment = -4152
    }
}

# Income statement build, column M (Q1'24).
Copy-ModelStyle 15
$wsModel.Range("M15").Value = 7541

Copy-ModelStyle 16
$wsModel.Range("M16").Value = 7025

Copy-ModelStyle 17
$wsModel.Range("M17").Value = 186

Copy-ModelStyle 18
$wsModel.Range("M18").Formula = "=SUM(M15:M17)"

Copy-ModelStyle 19
$wsModel.Range("M19").Value = 5217

Copy-ModelStyle 20
$wsModel.Range("M20").Value = 1419

Copy-ModelStyle 21
$wsModel.Range("M21").Value = 94

Copy-ModelStyle 22
$wsModel.Range("M22").Formula = "=SUM(M19:M21)"

Copy-ModelStyle 23
$wsModel.Range("M23").Formula = "=M18-M22"

Copy-ModelStyle 24
$wsModel.Range("M24").Value = 4458

Copy-ModelStyle 25
$wsModel.Range("M25").Value = 1685

Copy-ModelStyle 26
$wsModel.Range("M26").Formula = "=M24+M25"

Copy-ModelStyle 27
$wsModel.Range("M27").Formula = "=M23-M26"

# Row 28 (Interest Income) carries no explicit style in any quarter column.
$wsModel.Range("M28").Formula = "=-190-215+412"

Copy-ModelStyle 29
$wsModel.Range("M29").Formula = "=M27+M28"

# Row 30 (Taxes) also carries no explicit style.
$wsModel.Range("M30").Value = 159

Copy-ModelStyle 31
$wsModel.Range("M31").Formula = "=+M29-M30"

Copy-ModelStyle 32
$wsModel.Range("M32").Formula = "=M31/M33"

Copy-ModelStyle 33
$wsModel.Range("M33").Value = 912.8

Copy-ModelStyle 35
$wsModel.Range("M35").Formula = "=M18/I18-1"

Copy-ModelStyle 36
$wsModel.Range("M36").Formula = "=M23/M18"

# Q35's formula is untouched - once M18 is populated above, its cached
# #DIV/0! resolves to a real ratio automatically on recalculation.

# ---------------------------------------------------------------------------
# View state: Model keeps its own zoom/freeze, but the last touched cell
# becomes its selection; Main becomes the active tab/sheet with a new zoom
# and selection.
# ---------------------------------------------------------------------------
$wsModel.Activate()
$wsModel.Range("Q35").Select()

$wsMain.Activate()
$excel.ActiveWindow.Zoom = 180
$wsMain.Range("N7").Select()
